# Progress.xlsx update — "Add files via upload"
#
# БИВТ-22-17: student in row 22 gets a lab score (E22) instead of the
# "pass" placeholder, plus a new "failed 3.5" remark in F22.
#
# БИВТ-22-20: several students get lab scores / remarks filled in
# (rows 2, 8, 13, 16, 22, 24).
#
# Finally, the active sheet/selection bookmarks are moved around:
# БИВТ-22-17 -> L15, БИВТ-22-18 loses the tab-selected flag, and
# БИВТ-22-20 becomes the active sheet with selection E17.

$wb = $excel.ActiveWorkbook

# ---- БИВТ-22-20 ---------------------------------------------------------
$ws3 = $wb.Worksheets.Item("БИВТ-22-20")

$ws3.Range("F2").Value = 5

$ws3.Range("E8").Value = "failed 2.20"

# ---- БИВТ-22-17 ---------------------------------------------------------
$ws1 = $wb.Worksheets.Item("БИВТ-22-17")

$ws1.Range("E22").Value = 5
$ws1.Range("F22").Value = "failed 3.5"

# ---- back to БИВТ-22-20 --------------------------------------------------
$ws3.Range("F13").Value = 5
$ws3.Range("G13").Value = "pass"

$ws3.Range("E16").Value = "failed 2.18"

$ws3.Range("E22").Value = 5

$ws3.Range("E24").Value = 5

# ---- selections / active sheet ------------------------------------------
# Touch БИВТ-22-17 first so its own selection bookmark is updated ...
[void]$ws1.Range("L15").Select()

# ... then land on БИВТ-22-20 last so it ends up as the active tab.
[void]$ws3.Activate()
[void]$ws3.Range("E17").Select()
